$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.716.85'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '1.600.42'
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = '''211.13'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").Value = '''0.0844'
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("D12").Value = '1.825.40'
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").Value = '1.598.98'
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").Value = '''4.03'
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").Value = '''65.01'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").Value = '26.691.49'
$ws.Range("D18").Value = '0.0₃0739'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").Value = '''210.18'
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("D20").Value = '''7.20'
$ws.Range("E20").Value = '  +2.55%  '
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("D25").Value = '''144.04'
$ws.Range("E25").Value = '  -0.76%  '
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").Value = '''7.08'
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("D29").Value = '''15.34'
$ws.Range("E29").Value = '  +0.53%  '
$ws.Range("D30").Value = '''0.0512'
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("E32").Value = '  +1.11%  '
$ws.Range("E33").Value = '  +1.14%  '
$ws.Range("D34").Value = '1.292.54'
$ws.Range("E34").Value = '  +1.23%  '
$ws.Range("E35").Value = '  +0.61%  '
$ws.Range("E36").Value = '  +0.75%  '
$ws.Range("D37").Value = '''0.601'
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("E38").Value = '  +13.35%  '
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("E40").Value = '  -2.05%  '
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("D43").Value = '''0.779'
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("D44").Value = '''63.04'
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("D45").Value = '1.738.42'
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").Value = '''90.64'
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("E47").Value = '  -3.23%  '
$ws.Range("E48").Value = '  -0.86%  '
$ws.Range("E49").Value = '  +1.75%  '
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").Value = '''7.41'
$ws.Range("E51").Value = '  -0.62%  '
